$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.361143666666667
$ws.Range("H2").Value = 28.083431
$ws.Range("I2").Value = 0.1965934861218525
$ws.Range("J2").Value = 0.1965934861218526
$ws.Range("M2").Value = 10.82167433333333
$ws.Range("N2").Value = 32.465023
$ws.Range("O2").Value = 0.09133543757015983
$ws.Range("P2").Value = 0.09133543757015983
$ws.Range("Q2").Value = 101.3032481482126
$ws.Range("R2").Value = 911.729233333913
$ws.Range("S2").Value = 0.01795595207838254
$ws.Range("T2").Value = 0.01795595207838255

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.361143666666667
$ws.Range("H3").Value = 28.083431
$ws.Range("I3").Value = 0.1965934861218525
$ws.Range("J3").Value = 0.1965934861218526
$ws.Range("M3").Value = 36.14140700000001
$ws.Range("O3").Value = 0.3050351656377608
$ws.Range("P3").Value = 0.3050351656377608
$ws.Range("Q3").Value = 338.3249032424724
$ws.Range("R3").Value = 3044.924129182251
$ws.Range("S3").Value = 0.05996792660248412
$ws.Range("T3").Value = 0.05996792660248412

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.361143666666667
$ws.Range("H4").Value = 28.083431
$ws.Range("I4").Value = 0.1965934861218525
$ws.Range("J4").Value = 0.1965934861218526
$ws.Range("M4").Value = 26.40107466666666
$ws.Range("N4").Value = 79.20322399999999
$ws.Range("O4").Value = 0.2228263051286729
$ws.Range("P4").Value = 0.2228263051286729
$ws.Range("Q4").Value = 247.1442529090604
$ws.Range("R4").Value = 2224.298276181544
$ws.Range("S4").Value = 0.04380620012489743
$ws.Range("T4").Value = 0.04380620012489744

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.361143666666667
$ws.Range("H5").Value = 28.083431
$ws.Range("I5").Value = 0.1965934861218525
$ws.Range("J5").Value = 0.1965934861218526
$ws.Range("M5").Value = 45.11859966666666
$ws.Range("N5").Value = 135.355799
$ws.Range("O5").Value = 0.3808030916634065
$ws.Range("P5").Value = 0.3808030916634065
$ws.Range("Q5").Value = 422.3616935184854
$ws.Range("R5").Value = 3801.255241666369
$ws.Range("S5").Value = 0.07486340731608844
$ws.Range("T5").Value = 0.07486340731608847

$ws.Range("I6").Value = 0.2809659460057553
$ws.Range("J6").Value = 0.2809659460057554
$ws.Range("M6").Value = 10.82167433333333
$ws.Range("N6").Value = 32.465023
$ws.Range("O6").Value = 0.09133543757015983
$ws.Range("P6").Value = 0.09133543757015983
$ws.Range("Q6").Value = 144.7797865071508
$ws.Range("R6").Value = 1303.018078564357
$ws.Range("S6").Value = 0.02566214762074956
$ws.Range("T6").Value = 0.02566214762074957

$ws.Range("I7").Value = 0.2809659460057553
$ws.Range("J7").Value = 0.2809659460057554
$ws.Range("M7").Value = 36.14140700000001
$ws.Range("O7").Value = 0.3050351656377608
$ws.Range("P7").Value = 0.3050351656377608
$ws.Range("Q7").Value = 483.5245478983378
$ws.Range("R7").Value = 4351.72093108504
$ws.Range("S7").Value = 0.08570449387843572
$ws.Range("T7").Value = 0.08570449387843573

$ws.Range("I8").Value = 0.2809659460057553
$ws.Range("J8").Value = 0.2809659460057554
$ws.Range("M8").Value = 26.40107466666666
$ws.Range("N8").Value = 79.20322399999999
$ws.Range("O8").Value = 0.2228263051286729
$ws.Range("P8").Value = 0.2228263051286729
$ws.Range("Q8").Value = 353.2116968282462
$ws.Range("R8").Value = 3178.905271454216
$ws.Range("S8").Value = 0.06260660361544466
$ws.Range("T8").Value = 0.06260660361544468

$ws.Range("I9").Value = 0.2809659460057553
$ws.Range("J9").Value = 0.2809659460057554
$ws.Range("M9").Value = 45.11859966666666
$ws.Range("N9").Value = 135.355799
$ws.Range("O9").Value = 0.3808030916634065
$ws.Range("P9").Value = 0.3808030916634065
$ws.Range("Q9").Value = 603.6275927395712
$ws.Range("R9").Value = 5432.648334656141
$ws.Range("S9").Value = 0.1069927008911254
$ws.Range("T9").Value = 0.1069927008911254

$ws.Range("G10").Value = 7.684952333333334
$ws.Range("H10").Value = 23.054857
$ws.Range("I10").Value = 0.1613917725961189
$ws.Range("J10").Value = 0.1613917725961189
$ws.Range("M10").Value = 10.82167433333333
$ws.Range("N10").Value = 32.465023
$ws.Range("O10").Value = 0.09133543757015983
$ws.Range("P10").Value = 0.09133543757015983
$ws.Range("Q10").Value = 83.16405141852346
$ws.Range("R10").Value = 748.4764627667112
$ws.Range("S10").Value = 0.01474078817029025
$ws.Range("T10").Value = 0.01474078817029025

$ws.Range("G11").Value = 7.684952333333334
$ws.Range("H11").Value = 23.054857
$ws.Range("I11").Value = 0.1613917725961189
$ws.Range("J11").Value = 0.1613917725961189
$ws.Range("M11").Value = 36.14140700000001
$ws.Range("O11").Value = 0.3050351656377608
$ws.Range("P11").Value = 0.3050351656377608
$ws.Range("Q11").Value = 277.7449900545998
$ws.Range("R11").Value = 2499.704910491398
$ws.Range("S11").Value = 0.04923016608642895
$ws.Range("T11").Value = 0.04923016608642895

$ws.Range("G12").Value = 7.684952333333334
$ws.Range("H12").Value = 23.054857
$ws.Range("I12").Value = 0.1613917725961189
$ws.Range("J12").Value = 0.1613917725961189
$ws.Range("M12").Value = 26.40107466666666
$ws.Range("N12").Value = 79.20322399999999
$ws.Range("O12").Value = 0.2228263051286729
$ws.Range("P12").Value = 0.2228263051286729
$ws.Range("Q12").Value = 202.8910003621075
$ws.Range("R12").Value = 1826.019003258968
$ws.Range("S12").Value = 0.03596233236576017
$ws.Range("T12").Value = 0.03596233236576018

$ws.Range("G13").Value = 7.684952333333334
$ws.Range("H13").Value = 23.054857
$ws.Range("I13").Value = 0.1613917725961189
$ws.Range("J13").Value = 0.1613917725961189
$ws.Range("M13").Value = 45.11859966666666
$ws.Range("N13").Value = 135.355799
$ws.Range("O13").Value = 0.3808030916634065
$ws.Range("P13").Value = 0.3808030916634065
$ws.Range("Q13").Value = 346.7342877850825
$ws.Range("R13").Value = 3120.608590065743
$ws.Range("S13").Value = 0.06145848597363952
$ws.Range("T13").Value = 0.06145848597363952

$ws.Range("G14").Value = 17.19197166666666
$ws.Range("H14").Value = 51.57591499999999
$ws.Range("I14").Value = 0.3610487952762732
$ws.Range("J14").Value = 0.3610487952762733
$ws.Range("M14").Value = 10.82167433333333
$ws.Range("N14").Value = 32.465023
$ws.Range("O14").Value = 0.09133543757015983
$ws.Range("P14").Value = 0.09133543757015983
$ws.Range("Q14").Value = 186.0459185245605
$ws.Range("R14").Value = 1674.413266721045
$ws.Range("S14").Value = 0.03297654970073747
$ws.Range("T14").Value = 0.03297654970073748

$ws.Range("G15").Value = 17.19197166666666
$ws.Range("H15").Value = 51.57591499999999
$ws.Range("I15").Value = 0.3610487952762732
$ws.Range("J15").Value = 0.3610487952762733
$ws.Range("M15").Value = 36.14140700000001
$ws.Range("O15").Value = 0.3050351656377608
$ws.Range("P15").Value = 0.3050351656377608
$ws.Range("Q15").Value = 621.3420451374684
$ws.Range("R15").Value = 5592.078406237215
$ws.Range("S15").Value = 0.110132579070412
$ws.Range("T15").Value = 0.110132579070412

$ws.Range("G16").Value = 17.19197166666666
$ws.Range("H16").Value = 51.57591499999999
$ws.Range("I16").Value = 0.3610487952762732
$ws.Range("J16").Value = 0.3610487952762733
$ws.Range("M16").Value = 26.40107466666666
$ws.Range("N16").Value = 79.20322399999999
$ws.Range("O16").Value = 0.2228263051286729
$ws.Range("P16").Value = 0.2228263051286729
$ws.Range("Q16").Value = 453.8865276388843
$ws.Range("R16").Value = 4084.978748749959
$ws.Range("S16").Value = 0.0804511690225706
$ws.Range("T16").Value = 0.08045116902257063

$ws.Range("G17").Value = 17.19197166666666
$ws.Range("H17").Value = 51.57591499999999
$ws.Range("I17").Value = 0.3610487952762732
$ws.Range("J17").Value = 0.3610487952762733
$ws.Range("M17").Value = 45.11859966666666
$ws.Range("N17").Value = 135.355799
$ws.Range("O17").Value = 0.3808030916634065
$ws.Range("P17").Value = 0.3808030916634065
$ws.Range("Q17").Value = 775.6776871090092
$ws.Range("R17").Value = 6981.099183981084
$ws.Range("S17").Value = 0.1374884974825532
$ws.Range("T17").Value = 0.1374884974825532
